# Hortaliza, Vega Modelo de Temuco - Puerro: add a new weekly price record.
# A new row of data is inserted at row 203, pushing the existing rows
# (previously 203-248) down to 204-249, growing the sheet from 248 to 249 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 203, shifting rows 203:248 -> 204:249
$ws.Rows.Item(203).Insert()

# Populate the new row 203 with the new weekly record
$ws.Range("A203").Value = 10
$ws.Range("B203").Value = "Vega Modelo de Temuco"
$ws.Range("C203").Value = "La Araucanía"
$ws.Range("D203").Value = 44855
$ws.Range("E203").Value = 9
$ws.Range("F203").Value = 100112005
$ws.Range("G203").Value = "Puerro"
$ws.Range("H203").Value = "Azul de Maquehue"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 20
$ws.Range("K203").Value = 1600
$ws.Range("L203").Value = 1600
$ws.Range("M203").Value = 1600
$ws.Range("N203").Value = "`$/docena de paquetes"
$ws.Range("O203").Value = "Provincia de Cautín"
$ws.Range("P203").Value = 133
$ws.Range("Q203").Value = 12
$ws.Range("R203").Value = "Hortaliza"
